$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 7937.75  # H76
$ws.Cells.Item(76, 10).Value = 8071.2856  # J76
$ws.Cells.Item(76, 12).Value = 8071.2856  # L76
$ws.Cells.Item(76, 14).Value = -8701.285599999999  # N76
$ws.Cells.Item(79, 8).Value = 7937.75  # H79
$ws.Cells.Item(79, 10).Value = 8071.2856  # J79
$ws.Cells.Item(79, 12).Value = 8071.2856  # L79
$ws.Cells.Item(79, 14).Value = -10255.2856  # N79
$ws.Cells.Item(103, 8).Value = 2299.4  # H103
$ws.Cells.Item(103, 9).Value = 2745.5  # I103
$ws.Cells.Item(103, 10).Value = 2230.7693  # J103
$ws.Cells.Item(103, 11).Value = 8236.5  # K103
$ws.Cells.Item(103, 12).Value = 6692.3079  # L103
$ws.Cells.Item(103, 13).Value = -7650.5  # M103
$ws.Cells.Item(103, 14).Value = -7864.3079  # N103
$ws.Cells.Item(111, 8).Value = 1333.3334  # H111
$ws.Cells.Item(111, 9).Value = 625  # I111
$ws.Cells.Item(111, 11).Value = 1875  # K111
$ws.Cells.Item(111, 13).Value = 1192  # M111
$ws.Cells.Item(137, 8).Value = 19998  # H137
$ws.Cells.Item(137, 9).Value = 19998  # I137
$ws.Cells.Item(137, 10).Value = 0  # J137
$ws.Cells.Item(137, 11).Value = 59994  # K137
$ws.Cells.Item(137, 12).Value = 0  # L137
$ws.Cells.Item(137, 13).Value = -57444  # M137
$ws.Cells.Item(137, 14).ClearContents()  # N137 was -12600

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 4235.8237  # H63
$ws.Cells.Item(63, 9).Value = 1139.8889  # I63
$ws.Cells.Item(63, 10).Value = 7718.75  # J63
$ws.Cells.Item(63, 11).Value = 1139.8889  # K63
$ws.Cells.Item(63, 12).Value = 7718.75  # L63
$ws.Cells.Item(63, 13).Value = -453.8888999999999  # M63
$ws.Cells.Item(63, 14).Value = -9090.75  # N63
$ws.Cells.Item(66, 8).Value = 4235.8237  # H66
$ws.Cells.Item(66, 9).Value = 1139.8889  # I66
$ws.Cells.Item(66, 10).Value = 7718.75  # J66
$ws.Cells.Item(66, 11).Value = 5699.4445  # K66
$ws.Cells.Item(66, 12).Value = 38593.75  # L66
$ws.Cells.Item(66, 13).Value = -2267.4445  # M66
$ws.Cells.Item(66, 14).Value = -45457.75  # N66
$ws.Cells.Item(74, 8).Value = 3309.875  # H74
$ws.Cells.Item(74, 9).Value = 3163.1667  # I74
$ws.Cells.Item(74, 11).Value = 3163.1667  # K74
$ws.Cells.Item(74, 13).Value = -2289.1667  # M74
$ws.Cells.Item(77, 8).Value = 3309.875  # H77
$ws.Cells.Item(77, 9).Value = 3163.1667  # I77
$ws.Cells.Item(77, 11).Value = 15815.8335  # K77
$ws.Cells.Item(77, 13).Value = -11447.8335  # M77
$ws.Cells.Item(132, 8).Value = 3000  # H132
$ws.Cells.Item(132, 9).Value = 3000  # I132
$ws.Cells.Item(132, 10).Value = 3000  # J132
$ws.Cells.Item(132, 11).Value = 9000  # K132
$ws.Cells.Item(132, 12).Value = 9000  # L132
$ws.Cells.Item(132, 13).Value = -6470  # M132
$ws.Cells.Item(132, 14).Value = -14060  # N132

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2278.611  # H134
$ws.Cells.Item(134, 9).Value = 2047.8125  # I134
$ws.Cells.Item(134, 11).Value = 6143.4375  # K134
$ws.Cells.Item(134, 13).Value = -3608.4375  # M134

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 1230.8  # H2
$ws.Cells.Item(2, 10).Value = 75  # J2
$ws.Cells.Item(2, 12).Value = 75  # L2
$ws.Cells.Item(2, 14).Value = -301  # N2
$ws.Cells.Item(95, 8).Value = 16021.1  # H95
$ws.Cells.Item(95, 10).Value = 16021.1  # J95
$ws.Cells.Item(95, 12).Value = 16021.1  # L95
$ws.Cells.Item(95, 14).Value = -21513.1  # N95

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(28, 8).Value = 2066.6667  # H28
$ws.Cells.Item(28, 9).Value = 1000  # I28
$ws.Cells.Item(28, 10).Value = 2600  # J28
$ws.Cells.Item(28, 11).Value = 3000  # K28
$ws.Cells.Item(28, 12).Value = 7800  # L28
$ws.Cells.Item(28, 13).Value = -2768  # M28
$ws.Cells.Item(28, 14).Value = -8264  # N28
$ws.Cells.Item(86, 8).Value = 81.25  # H86
$ws.Cells.Item(86, 10).Value = 78.333336  # J86
$ws.Cells.Item(86, 12).Value = 235.000008  # L86
$ws.Cells.Item(86, 14).Value = -2607.000008  # N86
$ws.Cells.Item(89, 8).Value = 81.25  # H89
$ws.Cells.Item(89, 10).Value = 78.333336  # J89
$ws.Cells.Item(89, 12).Value = 705.0000240000001  # L89
$ws.Cells.Item(89, 14).Value = -12561.000024  # N89
$ws.Cells.Item(104, 8).Value = 2000  # H104
$ws.Cells.Item(104, 9).Value = 2000  # I104
$ws.Cells.Item(104, 11).Value = 6000  # K104
$ws.Cells.Item(104, 13).Value = -3379  # M104
$ws.Cells.Item(120, 8).Value = 2785.8  # H120
$ws.Cells.Item(120, 9).Value = 2785.8  # I120
$ws.Cells.Item(120, 11).Value = 8357.400000000001  # K120
$ws.Cells.Item(120, 13).Value = -3519.400000000001  # M120
$ws.Cells.Item(125, 8).Value = 3000  # H125
$ws.Cells.Item(125, 10).Value = 3000  # J125
$ws.Cells.Item(125, 12).Value = 9000  # L125
$ws.Cells.Item(125, 14).Value = -18840  # N125
$ws.Cells.Item(129, 8).Value = 1211.3334  # H129
$ws.Cells.Item(129, 9).Value = 199.5  # I129
$ws.Cells.Item(129, 10).Value = 1717.25  # J129
$ws.Cells.Item(129, 11).Value = 598.5  # K129
$ws.Cells.Item(129, 12).Value = 5151.75  # L129
$ws.Cells.Item(129, 13).Value = 4401.5  # M129
$ws.Cells.Item(129, 14).Value = -15151.75  # N129
$ws.Cells.Item(131, 8).Value = 2485.14  # H131
$ws.Cells.Item(131, 9).Value = 937  # I131
$ws.Cells.Item(131, 10).Value = 2549.6458  # J131
$ws.Cells.Item(131, 11).Value = 2811  # K131
$ws.Cells.Item(131, 12).Value = 7648.937399999999  # L131
$ws.Cells.Item(131, 13).Value = 2229  # M131
$ws.Cells.Item(131, 14).Value = -17728.9374  # N131
$ws.Cells.Item(137, 8).Value = 500  # H137
$ws.Cells.Item(137, 10).Value = 0  # J137
$ws.Cells.Item(137, 12).Value = 0  # L137
$ws.Cells.Item(137, 14).ClearContents()  # N137 was -13350
$ws.Cells.Item(138, 8).Value = 2520.8948  # H138
$ws.Cells.Item(138, 9).Value = 2632.3333  # I138
$ws.Cells.Item(138, 10).Value = 2500  # J138
$ws.Cells.Item(138, 11).Value = 7896.999899999999  # K138
$ws.Cells.Item(138, 12).Value = 7500  # L138
$ws.Cells.Item(138, 13).Value = -2756.999899999999  # M138
$ws.Cells.Item(138, 14).Value = -17780  # N138

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 34368.273  # H20
$ws.Cells.Item(20, 9).Value = 500  # I20
$ws.Cells.Item(20, 10).Value = 37755.1  # J20
$ws.Cells.Item(20, 11).Value = 500  # K20
$ws.Cells.Item(20, 12).Value = 37755.1  # L20
$ws.Cells.Item(20, 13).Value = -255  # M20
$ws.Cells.Item(20, 14).Value = -38245.1  # N20
$ws.Cells.Item(24, 8).Value = 13379.697  # H24
$ws.Cells.Item(24, 9).Value = 4532.6665  # I24
$ws.Cells.Item(24, 10).Value = 28862  # J24
$ws.Cells.Item(24, 11).Value = 4532.6665  # K24
$ws.Cells.Item(24, 12).Value = 28862  # L24
$ws.Cells.Item(24, 13).Value = -4359.6665  # M24
$ws.Cells.Item(24, 14).Value = -29208  # N24
$ws.Cells.Item(53, 8).Value = 29910.5  # H53
$ws.Cells.Item(53, 9).Value = 29999  # I53
$ws.Cells.Item(53, 10).Value = 29881  # J53
$ws.Cells.Item(53, 11).Value = 29999  # K53
$ws.Cells.Item(53, 12).Value = 29881  # L53
$ws.Cells.Item(53, 13).Value = -29368  # M53
$ws.Cells.Item(53, 14).Value = -31143  # N53
$ws.Cells.Item(122, 8).Value = 40145.258  # H122
$ws.Cells.Item(122, 9).Value = 2502.2222  # I122
$ws.Cells.Item(122, 11).Value = 7506.6666  # K122
$ws.Cells.Item(122, 13).Value = -5056.6666  # M122
$ws.Cells.Item(126, 8).Value = 4023.5386  # H126
$ws.Cells.Item(126, 9).Value = 3279.8  # I126
$ws.Cells.Item(126, 11).Value = 9839.400000000001  # K126
$ws.Cells.Item(126, 13).Value = -7369.400000000001  # M126

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2086.5  # H40
$ws.Cells.Item(40, 9).Value = 2092.0908  # I40
$ws.Cells.Item(40, 10).Value = 2066  # J40
$ws.Cells.Item(40, 11).Value = 2092.0908  # K40
$ws.Cells.Item(40, 12).Value = 2066  # L40
$ws.Cells.Item(40, 13).Value = -1956.0908  # M40
$ws.Cells.Item(40, 14).Value = -2338  # N40
$ws.Cells.Item(88, 8).Value = 22000  # H88
$ws.Cells.Item(88, 9).Value = 22000  # I88
$ws.Cells.Item(88, 11).Value = 22000  # K88
$ws.Cells.Item(88, 13).Value = -21572  # M88
$ws.Cells.Item(91, 8).Value = 22000  # H91
$ws.Cells.Item(91, 9).Value = 22000  # I91
$ws.Cells.Item(91, 11).Value = 22000  # K91
$ws.Cells.Item(91, 13).Value = -20518  # M91
$ws.Cells.Item(100, 8).Value = 8571.571  # H100
$ws.Cells.Item(100, 9).Value = 5833.8335  # I100
$ws.Cells.Item(100, 11).Value = 5833.8335  # K100
$ws.Cells.Item(100, 13).Value = -5292.8335  # M100

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 11704.833  # H14
$ws.Cells.Item(14, 9).Value = 5158  # I14
$ws.Cells.Item(14, 10).Value = 13887.111  # J14
$ws.Cells.Item(14, 11).Value = 5158  # K14
$ws.Cells.Item(14, 12).Value = 13887.111  # L14
$ws.Cells.Item(14, 13).Value = -4990  # M14
$ws.Cells.Item(14, 14).Value = -14223.111  # N14
